$wb = $excel.ActiveWorkbook

# Helper: assign a value that should be stored as text (not auto-converted
# to a number) while leaving the cell's number format back at General/Normal
# once written, so styles.xml stays effectively unchanged.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Sheet 1: GEOSPATIAL_FORMAT_NOT_MAPPED ---
# Insert a new data row (row 2) before the existing one, then fill both rows.
$ws1 = $wb.Worksheets.Item("GEOSPATIAL_FORMAT_NOT_MAPPED")
$ws1.Rows.Item(2).Insert()

$ws1.Range("A2").Value = "30aeb5c1-4285-46c8-b60b-15b1a6f4258b"
$ws1.Range("B2").Value = "fgdb"
Set-TextValue $ws1.Range("C2") "20201113"

$ws1.Range("A3").Value = "b6b86630-26c0-4b6d-a4b0-6633d7f00c7a"
$ws1.Range("B3").Value = "toto"
Set-TextValue $ws1.Range("C3") "20201113"

# --- Sheet 2: SPATIAL_TYPE_NOT_MAPPED ---
# Remove the Fort Sheppard data row (row 3), keep row 2 but refresh its timestamp.
$ws2 = $wb.Worksheets.Item("SPATIAL_TYPE_NOT_MAPPED")
$ws2.Rows.Item(3).Delete()
Set-TextValue $ws2.Range("D2") "20201113"

# --- Sheet 3: MAP_SERVICE_URL_ERROR ---
# Only the timestamp changes.
$ws3 = $wb.Worksheets.Item("MAP_SERVICE_URL_ERROR")
Set-TextValue $ws3.Range("E2") "20201113"

# --- Sheet 4: DELETED_RECORDS ---
# Not directly touched by the authored change, but every _timestamp cell here
# shares the same "20200915" shared string used on sheets 1-3. That shared
# string is being fully replaced with "20201113" (no cell anywhere keeps the
# old value), so refresh these cells too to keep the shared string table
# consistent with the rest of the workbook.
$ws4 = $wb.Worksheets.Item("DELETED_RECORDS")
for ($r = 2; $r -le 28; $r++) {
    Set-TextValue $ws4.Cells.Item($r, 2) "20201113"
}
